$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "datatype" column (K) describing the socket-stream field type ---
$ws.Range("K1").Value = "datatype"
$ws.Range("K1").Font.Bold = $true

$ws.Range("K2").Value  = "date"
$ws.Range("K3").Value  = "number"
$ws.Range("K4").Value  = "number"
$ws.Range("K5").Value  = "number"
$ws.Range("K6").Value  = "number"
$ws.Range("K7").Value  = "number"
$ws.Range("K8").Value  = "number"
$ws.Range("K9").Value  = "number"
$ws.Range("K10").Value = "number"
$ws.Range("K11").Value = "string"

# --- timestamp row now expressed as Julian dates instead of raw seconds ---
$ws.Range("D2").Value = 2415020.5
$ws.Range("E2").Value = 2444239.5
$ws.Range("G2").Value = 2469807.5
$ws.Range("H2").Value = 2524593.5
$ws.Range("G2:H2").NumberFormat = "0.00;[Red]0.00"

$ws.Range("J2").Value = "Julian date"

# --- cosmetics: widen warn_high column now that it holds Julian dates, and
#     move the active selection to reflect where the user was working ---
$ws.Columns("G").ColumnWidth = 11
$ws.Range("D4").Select()
